$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values regenerated to use K instead of Strike# for rows 2-9.
$kValues = @{
    2 = 1
    3 = 5
    4 = 1
    5 = 5
    6 = 1
    7 = 4
    8 = 3
    9 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
